# IM Asset UI changes
#
# The deck's "ADW+DVD for ALL LOBs - Introduction" table (Table 19) lists,
# for each Line of Business, the assets available (Video / Demo Walkthrough /
# DVD Live Demo / Demo Data Sets / Video Script). The header-row "Video"
# cell used to carry a second line with the clip's running time ("3:26").
# That duration line is removed, leaving just the "Video" label.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$tableShape = $s.Shapes.Item("Table 19")
$tbl = $tableShape.Table

# Row 1 ("Line of Business" header row), Column 2 ("Video" / "3:26").
$cell = $tbl.Cell(1, 2)
$cell.Shape.TextFrame.TextRange.Text = "Video"

# --- Incidental, auto-generated metadata ---------------------------------
# Re-saving the deck in PowerPoint after this edit also refreshed the two
# "Update automatically" date placeholders living in the Handout Master and
# the Notes Master (28-05-2019 -> 03-06-2019), and bumped the table's
# internal p14:modId co-authoring stamp. Those are side effects of
# PowerPoint's own save pipeline rather than explicit user actions, and the
# automation surface here has no field/modId writer, so these are
# best-effort / non-fatal if unsupported.
try {
    $hf = $p.HandoutMaster.HeadersFooters.DateAndTime
    $hf.Value = "03-06-2019"
} catch {
}

try {
    $hf2 = $p.NotesMaster.HeadersFooters.DateAndTime
    $hf2.Value = "03-06-2019"
} catch {
}
